$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("B2").Value = 0.55

$ws.Range("B3").Value = 0.4846153846153846
$ws.Range("C3").Value = 0.4576923076923077

$ws.Range("C4").Value = 0.6

# Row 5: rename c-FUDGE -> WildguardCTG, update C5
$ws.Range("A5").Value = "WildguardCTG"
$ws.Range("C5").Value = 0.6115384615384616

# New row 6: c-FUDGE with its new values, copying style from A5
$ws.Range("A6").Value = "c-FUDGE"
$ws.Range("B6").Value = 0.5615384615384615
$ws.Range("C6").Value = 0.5384615384615384

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
